$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (Q5)
$ws.Range("B7").Value = 0.05059694089399255
$ws.Range("C7").Value = 0.3754832313001323
$ws.Range("D7").Value = 0.2169190185971216
$ws.Range("E7").Value = 0.4657456586991677
$ws.Range("F7").Value = 0.4910741687265304
$ws.Range("G7").Value = 9

# Row 8 (Q6)
$ws.Range("B8").Value = 0.1028329319255051
$ws.Range("C8").Value = 0.3256379015226427
$ws.Range("D8").Value = 0.1450293133628471
$ws.Range("E8").Value = 0.3808271436791856
$ws.Range("F8").Value = 0.3889235646740345
$ws.Range("G8").Value = 9
